$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{ row = 2;  B = -0.3684425897240435; C = 1.990300204684291;  D = 12.34821194374062; E = 3.514002268602088; F = 3.576871232156456; G = 22 }
    @{ row = 3;  B = -0.4406390919636032; C = 1.050095607268688;  D = 4.94590582631831;  E = 2.223939258684533; F = 2.233680856299617; G = 21 }
    @{ row = 4;  B = -0.2266117036152265; C = 0.7762331526076556; D = 2.687494923465916;  E = 1.639358082746389; F = 1.665799020938133; G = 20 }
    @{ row = 5;  B = 0.02747715458826221; C = 0.7115633747814063; D = 1.229190786505692;  E = 1.10868876899953;  F = 1.138719555200954; G = 19 }
    @{ row = 6;  B = -0.06432496402825524;C = 0.6672897034160007; D = 1.122755052931566;  E = 1.059601365104616; F = 1.088309869768877; G = 18 }
    @{ row = 7;  B = 0.1299503236878231;  C = 0.5419393736802232; D = 0.7287712587988439; E = 0.8536810052934549;F = 0.8696993347893518;G = 17 }
    @{ row = 8;  B = 0.1289229671324676;  C = 0.5250093336297077; D = 0.4739960288771982; E = 0.6884736951236395;F = 0.6984744499630291;G = 16 }
    @{ row = 9;  B = 0.196164091590951;   C = 0.5034467950485412; D = 0.4366068896125055; E = 0.6607623548693626;F = 0.6531187860533945;G = 15 }
    @{ row = 10; B = 0.2421836099529805;  C = 0.4042897487511507; D = 0.229010732524476;  E = 0.4785506582635489;F = 0.4283249702957092;G = 14 }
    @{ row = 11; B = 0.2725419345940574;  C = 0.3839238355923389; D = 0.1922556345917136; E = 0.4384696507076785;F = 0.3575023158704416;G = 13 }
)

foreach ($item in $data) {
    $r = $item.row
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 4).Value = $item.D
    $ws.Cells.Item($r, 5).Value = $item.E
    $ws.Cells.Item($r, 6).Value = $item.F
    $ws.Cells.Item($r, 7).Value = $item.G
}
